$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Enterprises density (per 1000 people) - SMEs / MSMEs columns
$ws.Range("C11").Value = "'5.81"
$ws.Range("D11").Value = "'17.81"

# Row 12: Employment (% of total) - Micro / SMEs columns
$ws.Range("B12").Value = "'14.31"
$ws.Range("C12").Value = "'35.69"

# Row 14: Enterprises (% of total) - Micro / SMEs / MSMEs columns
$ws.Range("B14").Value = "'63.02"
$ws.Range("C14").Value = "'30.53"
$ws.Range("D14").Value = "'93.55"
